# Apply odds updates to "Jogos da Semana" FlashScore sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 3.75
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("X2").Value = 9
$ws.Range("AS2").Value = 401

# Row 3 updates
$ws.Range("Q3").Value = 2.08
$ws.Range("R3").Value = 1.73

# Row 7 updates
$ws.Range("G7").Value = 2.55
$ws.Range("I7").Value = 2.5
$ws.Range("W7").Value = 9
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 15
$ws.Range("AJ7").Value = 26
$ws.Range("AK7").Value = 21
$ws.Range("AW7").Value = 4.75
$ws.Range("AX7").Value = 15
